# Religion.docx edits
# Applies the textual changes described by the commit:
#  - merges a couple of adjacent runs that had identical text/formatting
#    but were previously split (no visible text change)
#  - fixes "un de quatre éléments" -> "un des quatre éléments"
#  - fixes "que les reste du monde" -> "que le reste du monde" and moves
#    the _GoBack bookmark to sit right after "...que le"

$d = $word.ActiveDocument

# 1) "Dans les Cités Libres, ... Ostalyens). Parmi ces dieux figure "
#    was split across 3 runs (incl. a lone space run); re-merge them.
$d.Content.Find.Execute(
    "Dans les Cités Libres, d’autres dieux sont vénérés (même si certains vénèrent le Guerrier et la Vierge comme les Ostalyens). Parmi ces dieux figure ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Dans les Cités Libres, d’autres dieux sont vénérés (même si certains vénèrent le Guerrier et la Vierge comme les Ostalyens). Parmi ces dieux figure ",
    2
) | Out-Null

# 2) " (animaux officiellement, mais également des humains officieusement)."
#    was split across 2 runs; re-merge them.
$d.Content.Find.Execute(
    " (animaux officiellement, mais également des humains officieusement).",
    $false, $false, $false, $false, $false, $true, 1, $false,
    " (animaux officiellement, mais également des humains officieusement).",
    2
) | Out-Null

# 3) grammar fix: "un de quatre éléments" -> "un des quatre éléments"
#    (keep the original non-breaking space before the colon: "éléments : ")
$d.Content.Find.Execute(
    "chacune un de quatre éléments : ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "chacune un des quatre éléments : ",
    2
) | Out-Null

# 4) grammar fix: "que les reste du monde" -> "que le reste du monde"
$d.Content.Find.Execute(
    "Ostalya que les reste du monde ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Ostalya que le reste du monde ",
    2
) | Out-Null

# Move the "_GoBack" bookmark from after "...aussi quand" to right after
# "...Ostalya que le" (before " reste du monde ...").
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$target = $d.Content
$target.Find.Execute("Ostalya que le", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target.Collapse(0)
$d.Bookmarks.Add("_GoBack", $target) | Out-Null
